$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Sacrifice" action card, added as row 12.
$ws.Range("A12").Value = "action"
$ws.Range("B12").Value = "Sacrifice"
$ws.Range("C12").Value = "Sacrifice Self"
$ws.Range("D12").Value = "action_template.png"
$ws.Range("E12").Value = "blank.png"
$ws.Range("F12").Value = "-"
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = "The gladiator throws themselves in front of the opponent, saving their comrades."
$ws.Range("J12").Value = "The fighter is discarded`nTargeted creature cannot attack"

$ws.Range("J12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 30

$ws.Range("H13").Select()
